{"js": "const body = context.document.body;\n\nconst pairs = [\n  [\"2025-08-29 Friday\", \"2025-08-30 Saturday\"],\n  [\"20\u00d757=\", \"38\u00d781=\"],\n  [\"84\u00d781=\", \"60\u00d775=\"],\n  [\"79\u00d750=\", \"12\u00d734=\"],\n  [\"39\u00d776=\", \"90\u00d726=\"],\n  [\"68\u00d725=\", \"60\u00d744=\"],\n  [\"18\u00d722=\", \"92\u00d773=\"],\n  [\"21\u00d768=\", \"52\u00d713=\"],\n  [\"70\u00d761=\", \"30\u00d755=\"],\n  [\"89\u00d775=\", \"74\u00d752=\"],\n  [\"87\u00d732=\", \"53\u00d795=\"],\n  [\"68\u00d799=\", \"51\u00d752=\"],\n  [\"40\u00d780=\", \"44\u00d732=\"],\n  [\"60\u00d774=\", \"57\u00d711=\"],\n  [\"40\u00d739=\", \"99\u00d740=\"],\n  [\"40\u00d772=\", \"46\u00d754=\"],\n  [\"54\u00d733=\", \"66\u00d733=\"],\n  [\"95\u00d763=\", \"26\u00d775=\"],\n  [\"83\u00d768=\", \"38\u00d772=\"],\n  [\"44\u00d758=\", \"77\u00d725=\"],\n  [\"23\u00d763=\", \"83\u00d763=\"],\n  [\"62\u00d743=\", \"78\u00d714=\"],\n  [\"47\u00d717=\", \"40\u00d785=\"],\n  [\"79\u00d769=\", \"84\u00d718=\"],\n  [\"72\u00d788=\", \"75\u00d741=\"],\n  [\"98\u00d761=\", \"41\u00d777=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-08-29 Friday\", \"2025-08-30 Saturday\"),\n    @(\"20\u00d757=\", \"38\u00d781=\"),\n    @(\"84\u00d781=\", \"60\u00d775=\"),\n    @(\"79\u00d750=\", \"12\u00d734=\"),\n    @(\"39\u00d776=\", \"90\u00d726=\"),\n    @(\"68\u00d725=\", \"60\u00d744=\"),\n    @(\"18\u00d722=\", \"92\u00d773=\"),\n    @(\"21\u00d768=\", \"52\u00d713=\"),\n    @(\"70\u00d761=\", \"30\u00d755=\"),\n    @(\"89\u00d775=\", \"74\u00d752=\"),\n    @(\"87\u00d732=\", \"53\u00d795=\"),\n    @(\"68\u00d799=\", \"51\u00d752=\"),\n    @(\"40\u00d780=\", \"44\u00d732=\"),\n    @(\"60\u00d774=\", \"57\u00d711=\"),\n    @(\"40\u00d739=\", \"99\u00d740=\"),\n    @(\"40\u00d772=\", \"46\u00d754=\"),\n    @(\"54\u00d733=\", \"66\u00d733=\"),\n    @(\"95\u00d763=\", \"26\u00d775=\"),\n    @(\"83\u00d768=\", \"38\u00d772=\"),\n    @(\"44\u00d758=\", \"77\u00d725=\"),\n    @(\"23\u00d763=\", \"83\u00d763=\"),\n    @(\"62\u00d743=\", \"78\u00d714=\"),\n    @(\"47\u00d717=\", \"40\u00d785=\"),\n    @(\"79\u00d769=\", \"84\u00d718=\"),\n    @(\"72\u00d788=\", \"75\u00d741=\"),\n    @(\"98\u00d761=\", \"41\u00d777=\"),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
